$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows 3 and 4 (continuing the two existing RC-filter calc blocks)
# ---------------------------------------------------------------------------

# Row 3 - left block (A:D) & right block (H:K)
$ws.Range("A3").Value = 0.0015
$ws.Range("B3").Value = 4700
$ws.Range("C3").Formula = "=1/(2*PI()*A3*1000*B3*10^-9)"
$ws.Range("D3").Formula = "=9.2*C3*B3*10^-9*1000"

$ws.Range("H3").Value = 4700
$ws.Range("I3").Value = 470
$ws.Range("J3").Formula = "=1/(2*PI()*H3*10^-9*I3)"
$ws.Range("K3").Formula = "=9.2*H3*10^-9*I3*1000"

# Row 4 - left block (A:D) & right block (H:K)
$ws.Range("A4").Value = 0.1
$ws.Range("B4").Value = 4700
$ws.Range("C4").Formula = "=1/(2*PI()*A4*1000*B4*10^-9)"
$ws.Range("D4").Formula = "=9.2*C4*B4*10^-9*1000"

$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 47000
$ws.Range("J4").Formula = "=1/(2*PI()*H4*10^-9*I4)"
$ws.Range("K4").Formula = "=9.2*H4*10^-9*I4*1000"

# ---------------------------------------------------------------------------
# Notes about pwm / cutoff frequency choices
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "for stock arduino pwm of 491 Hz, want cutoff freq around 1Hz"
$ws.Range("B12").Value = "for sped up arduino pwm of 31372 Hz, try cutoff freq around 75 Hz"

# ---------------------------------------------------------------------------
# Column widths for the newly-used columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.17
$ws.Columns.Item(3).ColumnWidth = 9.33
$ws.Columns.Item(4).ColumnWidth = 8.6
$ws.Columns.Item(9).ColumnWidth = 9.33
$ws.Columns.Item(11).ColumnWidth = 8.6

# ---------------------------------------------------------------------------
# Turn both data blocks into Excel Tables
#   - H1:K4 keeps the default name "Table1" (table1.xml)
#   - A1:D4 is created second (default "Table2") and renamed to "Table3"
#     (table2.xml), matching the workbook's table history.
# ---------------------------------------------------------------------------
$tbl1 = $ws.ListObjects.Add(1, $ws.Range("H1:K4"), $null, 1)

$tbl2 = $ws.ListObjects.Add(1, $ws.Range("A1:D4"), $null, 1)
$tbl2.Name = "Table3"

# ---------------------------------------------------------------------------
# Selection moves on to the next empty row, as happened when the author
# finished entering data.
# ---------------------------------------------------------------------------
$ws.Range("B13").Select()
